$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh NATMI LR-pair edge-weight columns (Spp1-Itga4) with recomputed TPM values.
# Columns: G/H = ligand avg/total expr, I/J = ligand specificity (avg/total),
#          M/N = receptor avg/total expr, O/P = receptor specificity (avg/total),
#          Q/R = edge avg/total weight, S/T = edge specificity (avg/total).

# Row 2
$ws.Range("G2").Value = 0.7476426666666667
$ws.Range("H2").Value = 2.242928
$ws.Range("I2").Value = 0.001581772089386036
$ws.Range("J2").Value = 0.001581772089386036
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 0.1061101823235556
$ws.Range("R2").Value = 0.9549916409120001
$ws.Range("S2").Value = 0.000004085810280752765
$ws.Range("T2").Value = 0.000004085810280752765

# Row 3
$ws.Range("G3").Value = 0.7476426666666667
$ws.Range("H3").Value = 2.242928
$ws.Range("I3").Value = 0.001581772089386036
$ws.Range("J3").Value = 0.001581772089386036
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 0.05749795770844444
$ws.Range("R3").Value = 0.517481619376
$ws.Range("S3").Value = 0.000002213979295701374
$ws.Range("T3").Value = 0.000002213979295701375

# Row 4
$ws.Range("G4").Value = 0.7476426666666667
$ws.Range("H4").Value = 2.242928
$ws.Range("I4").Value = 0.001581772089386036
$ws.Range("J4").Value = 0.001581772089386036
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 1.954331201882667
$ws.Range("R4").Value = 17.588980816944
$ws.Range("S4").Value = 0.00007525221747616863
$ws.Range("T4").Value = 0.00007525221747616863

# Row 5
$ws.Range("G5").Value = 0.7476426666666667
$ws.Range("H5").Value = 2.242928
$ws.Range("I5").Value = 0.001581772089386036
$ws.Range("J5").Value = 0.001581772089386036
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 38.96133582407289
$ws.Range("R5").Value = 350.652022416656
$ws.Range("S5").Value = 0.001500220082333413
$ws.Range("T5").Value = 0.001500220082333413

# Row 6
$ws.Range("I6").Value = 0.02590993131491687
$ws.Range("J6").Value = 0.02590993131491688
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 1.738118629267111
$ws.Range("R6").Value = 15.643067663404
$ws.Range("S6").Value = 0.00006692687552805163
$ws.Range("T6").Value = 0.00006692687552805163

# Row 7
$ws.Range("I7").Value = 0.02590993131491687
$ws.Range("J7").Value = 0.02590993131491688
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("Q7").Value = 0.9418348856768889
$ws.Range("R7").Value = 8.476513971092
$ws.Range("S7").Value = 0.00003626568699068175
$ws.Range("T7").Value = 0.00003626568699068176

# Row 8
$ws.Range("I8").Value = 0.02590993131491687
$ws.Range("J8").Value = 0.02590993131491688
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 32.01256840170534
$ws.Range("R8").Value = 288.113115615348
$ws.Range("S8").Value = 0.001232655323220125
$ws.Range("T8").Value = 0.001232655323220125

# Row 9
$ws.Range("I9").Value = 0.02590993131491687
$ws.Range("J9").Value = 0.02590993131491688
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 638.1991071362058
$ws.Range("R9").Value = 5743.791964225852
$ws.Range("S9").Value = 0.02457408342917802
$ws.Range("T9").Value = 0.02457408342917802

# Row 10
$ws.Range("G10").Value = 18.93023433333333
$ws.Range("H10").Value = 56.79070299999999
$ws.Range("I10").Value = 0.04005030430848061
$ws.Range("J10").Value = 0.04005030430848062
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 2.686698748070778
$ws.Range("R10").Value = 24.180288732637
$ws.Range("S10").Value = 0.0001034522901174611
$ws.Range("T10").Value = 0.0001034522901174611

# Row 11
$ws.Range("G11").Value = 18.93023433333333
$ws.Range("H11").Value = 56.79070299999999
$ws.Range("I11").Value = 0.04005030430848061
$ws.Range("J11").Value = 0.04005030430848062
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 1.455842291561222
$ws.Range("R11").Value = 13.102580624051
$ws.Range("S11").Value = 0.00005605772482680046
$ws.Range("T11").Value = 0.00005605772482680048

# Row 12
$ws.Range("G12").Value = 18.93023433333333
$ws.Range("H12").Value = 56.79070299999999
$ws.Range("I12").Value = 0.04005030430848061
$ws.Range("J12").Value = 0.04005030430848062
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 49.48346217522433
$ws.Range("R12").Value = 445.351159577019
$ws.Range("S12").Value = 0.001905378296931734
$ws.Range("T12").Value = 0.001905378296931735

# Row 13
$ws.Range("G13").Value = 18.93023433333333
$ws.Range("H13").Value = 56.79070299999999
$ws.Range("I13").Value = 0.04005030430848061
$ws.Range("J13").Value = 0.04005030430848062
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 986.4969590054534
$ws.Range("R13").Value = 8878.472631049081
$ws.Range("S13").Value = 0.03798541599660462
$ws.Range("T13").Value = 0.03798541599660462

# Row 14
$ws.Range("G14").Value = 440.7369333333333
$ws.Range("H14").Value = 1322.2108
$ws.Range("I14").Value = 0.9324579922872165
$ws.Range("J14").Value = 0.9324579922872166
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 62.55217691257778
$ws.Range("R14").Value = 562.9695922131999
$ws.Range("S14").Value = 0.002408593802370088
$ws.Range("T14").Value = 0.002408593802370088

# Row 15
$ws.Range("G15").Value = 440.7369333333333
$ws.Range("H15").Value = 1322.2108
$ws.Range("I15").Value = 0.9324579922872165
$ws.Range("J15").Value = 0.9324579922872166
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 33.89516768262222
$ws.Range("R15").Value = 305.0565091436
$ws.Range("S15").Value = 0.001305145477586775
$ws.Range("T15").Value = 0.001305145477586775

# Row 16
$ws.Range("G16").Value = 440.7369333333333
$ws.Range("H16").Value = 1322.2108
$ws.Range("I16").Value = 0.9324579922872165
$ws.Range("J16").Value = 0.9324579922872166
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 1152.082377100933
$ws.Range("R16").Value = 10368.7413939084
$ws.Range("S16").Value = 0.04436134136759579
$ws.Range("T16").Value = 0.0443613413675958

# Row 17
$ws.Range("G17").Value = 440.7369333333333
$ws.Range("H17").Value = 1322.2108
$ws.Range("I17").Value = 0.9324579922872165
$ws.Range("J17").Value = 0.9324579922872166
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 22967.78987511684
$ws.Range("R17").Value = 206710.1088760516
$ws.Range("S17").Value = 0.8843829116396639
$ws.Range("T17").Value = 0.8843829116396639
